$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D values that look like plain numbers must be forced to Text format
# so Excel stores the exact literal string (preserving trailing zeros / precision)
# instead of silently converting them to a floating point number.
$ws.Range("D2").Value = "48.746.57"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.636.27"
$ws.Range("E3").Value = "  +4.45%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "110.06"
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "321.75"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.44"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.75"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.21"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").Value = "3.043.86"
$ws.Range("E15").Value = "  +4.27%  "
$ws.Range("D16").Value = "2.644.95"
$ws.Range("E16").Value = "  +7.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.858"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "48.758.25"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.66"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.89"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").Value = "0.0₃0940"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.70"
$ws.Range("E23").Value = "  -5.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.71"
$ws.Range("E24").Value = "  -1.99%  "
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.15"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.04"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  -4.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.30"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.19"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0793"
$ws.Range("E37").Value = "  +7.41%  "
$ws.Range("E38").Value = "  +2.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  +7.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.84"
$ws.Range("E40").Value = "  +3.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.69"
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.111"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.16"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("E44").Value = "  +3.02%  "
$ws.Range("D45").Value = "2.072.46"
$ws.Range("E45").Value = "  +3.36%  "
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.08"
$ws.Range("E47").Value = "  +5.26%  "
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.92"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "58.78"
$ws.Range("E50").Value = "  +4.30%  "
$ws.Range("E51").Value = "  -0.26%  "
